# Junction_Flooding_482 update:
#  - replace the data rows (2-5) with a new dataset (new timestamps + values)
#  - remove the old last data row (row 6) -> used range shrinks to A1:AH5
#  - a handful of column widths narrow by 1 character (E, I, P, T, Z)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) New values for rows 2-5 (columns A..AH), read left-to-right.
# ---------------------------------------------------------------------------
$newData = @(
    "45163.50694444445,7.205,5.605,1.192,15.61,12.374,4.794,14.798,8.957000000000001,4.429,5.631,6.249,7.306,2.788,6.015,7.738,5.138,0.492,0.931,84.14,16.424,5.552,10.174,6.283,0.894,9.741,4.361,4.885,6.06,8.134,1.522,13.18,3.625,6.387",
    "45163.51388888889,3.564,2.73,0.518,7.99,6.085,2.264,11.81,4.424,2.449,2.653,3.151,3.72,1.496,3.008,3.866,2.695,0.178,0.461,38.472,8.539,2.776,5.289,3.238,0.424,7.088,2.153,2.557,3.083,4.01,0.707,11.415,1.859,3.17",
    "45163.52083333334,1.725,1.312,0.31,3.988,2.875,0.966,7.24,2.156,1.359,1.177,1.557,1.913,0.867,1.504,1.853,1.423,0.09,0.286,15.585,4.348,1.388,2.614,1.695,0.211,4.259,1.02,1.388,1.652,2.064,0.461,7.111,0.987,1.531",
    "45163.52777777778,6.09,4.58,0.39,13.47,10.78,4.49,16.41,7.43,3.6,4.75,5.37,5.88,1.86,4.89,6.67,4.2,0.1,0.36,67.02,13.51,4.51,8.83,4.97,0.65,8.720000000000001,3.83,3.74,4.41,5.95,0.34,14.88,2.72,5.48"
)

$startRow = 2
for ($r = 0; $r -lt $newData.Length; $r++) {
    $values = $newData[$r].Split(",")
    $rowIndex = $startRow + $r
    for ($c = 0; $c -lt $values.Length; $c++) {
        $ws.Cells.Item($rowIndex, $c + 1).Value = [double]$values[$c]
    }
}

# ---------------------------------------------------------------------------
# 2) Drop the old trailing data row (was row 6) -- shifts nothing else,
#    just shrinks the used range down to A1:AH5.
# ---------------------------------------------------------------------------
$ws.Rows.Item(6).Delete()

# ---------------------------------------------------------------------------
# 3) Narrow a handful of columns by 1 character (OOXML `width`, i.e.
#    ColumnWidth = width - 0.83): E(5), I(9), P(16), T(20), Z(26).
# ---------------------------------------------------------------------------
$ws.Columns.Item(5).ColumnWidth = 6.17
$ws.Columns.Item(9).ColumnWidth = 6.17
$ws.Columns.Item(16).ColumnWidth = 6.17
$ws.Columns.Item(20).ColumnWidth = 7.17
$ws.Columns.Item(26).ColumnWidth = 6.17
